$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("First")

# C1 previously held a date serial (47574 == 2030-04-01) formatted with the
# custom "yyyy-mm-dd h:mm:ss" numFmt. Replace it with a plain text string
# "2030-04-01" (no date parsing, no special number format).
$ws.Range("C1").NumberFormat = "@"
$ws.Range("C1").Value = "2030-04-01"
$ws.Range("C1").Style = "Normal"

# Row 5: new set of values.
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 44
$ws.Range("C5").Value = 37
$ws.Range("D5").Value = 44
$ws.Range("E5").Value = 20
$ws.Range("F5").Value = 49
$ws.Range("G5").Value = 17
$ws.Range("H5").Value = 7
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 22
